$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at X:Y (existing X onward shifts right by 2)
$ws.Range("X1:Y1").EntireColumn.Insert()

# New columns pick up the same width as column W (their left neighbour)
$ws.Columns("X:Y").ColumnWidth = $ws.Columns("W").ColumnWidth

# Populate the new header cells.
# "derivation_description" is written first (Y1) so it lands earlier in the
# shared-string table, then "derived_variable" (X1), matching the source order.
$ws.Range("Y1").Value = "derivation_description"
$ws.Range("X1").Value = "derived_variable"

# Refresh the AutoFilter so its range grows to cover the two new columns.
$ws.AutoFilterMode = $false
$ws.Range("A1:AJ37").AutoFilter()

# Update the hidden _FilterDatabase defined name to match the new range.
foreach ($n in $wb.Names) {
  if ($n.Name -eq "Collection_EC!_FilterDatabase") {
    $n.RefersTo = "=Collection_EC!`$A`$1:`$AJ`$37"
  }
}
